# Atualização de bases das ligas, do dia: 16-05-2024 às 23:38
#
# The source feed re-shuffled several match rows inside their same
# date-block (the sequential index in column A stays put; every other
# column - id, teams, odds, results, P&L - follows its original row to
# its new position). This script captures the "before" values for each
# affected block and rewrites them into their new row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (id) and E..AB (everything except Div/Date, which are
# identical across each block anyway) - expressed as 1-based column
# indices: B=2, E=5, F=6, ... AB=28.
$cols = @(2,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28)

function Apply-RowPermutation($ws, $cols, $perm) {
    # $perm maps destination row -> source row (values as they were
    # BEFORE this block's edit). Snapshot every source row first so
    # overlapping/cyclical permutations don't clobber data we still need.
    $snapshot = @{}
    foreach ($destRow in $perm.Keys) {
        $srcRow = $perm[$destRow]
        if (-not $snapshot.ContainsKey($srcRow)) {
            $rowData = @{}
            foreach ($c in $cols) {
                $rowData[$c] = $ws.Cells.Item($srcRow, $c).Value()
            }
            $snapshot[$srcRow] = $rowData
        }
    }

    foreach ($destRow in $perm.Keys) {
        $srcRow = $perm[$destRow]
        $rowData = $snapshot[$srcRow]
        foreach ($c in $cols) {
            $ws.Cells.Item($destRow, $c).Value = $rowData[$c]
        }
    }
}

# Block 1: rows 93-99 (match date 45199.6875) - 7-way shuffle.
$perm1 = @{93=98; 94=99; 95=94; 96=97; 97=93; 98=95; 99=96}
Apply-RowPermutation $ws $cols $perm1

# Block 2: rows 114-115 (match date 45238.82291666666) - swap.
$perm2 = @{114=115; 115=114}
Apply-RowPermutation $ws $cols $perm2

# Block 3: rows 157-158 (match date 45360.79166666666) - swap.
$perm3 = @{157=158; 158=157}
Apply-RowPermutation $ws $cols $perm3

# Block 4: rows 173-174 (match date 45371.89583333334) - swap.
$perm4 = @{173=174; 174=173}
Apply-RowPermutation $ws $cols $perm4

# Block 5: rows 204 & 206 (match date 45401.83333333334) - swap.
$perm5 = @{204=206; 206=204}
Apply-RowPermutation $ws $cols $perm5

"done"
